$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -21.945
$ws.Range("A14").Value = -21.659
$ws.Range("A21").Value = -19.937
$ws.Range("A23").Value = -20.317
$ws.Range("A25").Value = -21.775
$ws.Range("A26").Value = -21.319
$ws.Range("A29").Value = -21.301
$ws.Range("A53").Value = -22.013
$ws.Range("A57").Value = -22.247
$ws.Range("A59").Value = -22.5
$ws.Range("A69").Value = -21.602
$ws.Range("A79").Value = -21.167
$ws.Range("A83").Value = -22.038
$ws.Range("A91").Value = -21.533
$ws.Range("A93").Value = -21.383
